# 03.06 update sequense diagrams
#
# Adds a "Диаграмма" (Diagram) row to the bottom of the UC-1 ("Регистрация
# клиента") and UC-12 ("Регистрация пользователя") use-case sheets, mirroring
# the row already present on sheets such as UC-3 / UC-14 / UC-16. Each new
# row links to a sequence-diagram image that is rendered elsewhere from the
# markdown-style text stored in the cell.

$wb = $excel.ActiveWorkbook

# Remember whatever sheet/selection was active so we can restore it - we
# only want to touch the two target sheets' own selections.
$originalActive = $wb.ActiveSheet

# Use an existing "Диаграмма" row as a formatting template (bold/bordered
# label cell + wrapped value cell) so the new rows pick up the same cell
# styles already defined in the workbook instead of minting new ones.
$templateSheet = $wb.Worksheets.Item("UC-3")
$templateRange = $templateSheet.Range("A8:B8")

# --- UC-1: Регистрация клиента -> registration_client.svg -----------------
$uc1 = $wb.Worksheets.Item("UC-1")
[void]$templateRange.Copy()
[void]$uc1.Range("A8:B8").PasteSpecial(-4122)
$uc1.Range("A8").Value = "Диаграмма"
$uc1.Range("B8").Value = '![](../diagrams/out/registration_client.svg){ width="100" }'
[void]$uc1.Range("A8:B8").Select()

# --- UC-12: Регистрация пользователя -> registration_user.svg -------------
$uc12 = $wb.Worksheets.Item("UC-12")
[void]$templateRange.Copy()
[void]$uc12.Range("A8:B8").PasteSpecial(-4122)
$uc12.Range("A8").Value = "Диаграмма"
$uc12.Range("B8").Value = '![](../diagrams/out/registration_user.svg){ width="100" }'
[void]$uc12.Range("A1:B8").Select()

# Restore original active sheet/selection.
[void]$originalActive.Activate()
